# "logea y vuelve a agendar con los parámetros correctos"
# Re-book reservations 2-5 (rows 2..5) as confirmed (estado = 1), then
# reselect the cell the user ends up on (F3), and let Excel's normal
# automatic recalculation refresh the TODAY()-based "fecha" formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# estado (column H) -> 1 for reservations in rows 2-5
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1

# Force the volatile TODAY()-based "fecha" formulas to recompute against
# the current clock so the cached values stay correct.
$ws.Calculate()

# Leave the selection on F3, matching where the user ended up after
# re-agendar-ing the second reservation.
[void]$ws.Range("F3").Select()
